# Remove the stray "ENTRY 2" markers that were added to the document.
$d = $word.ActiveDocument

# 1) "PRRB Decision in PRRB Case # 17-1920 ENTRY 2" -> "PRRB Decision in PRRB Case # 17-1920"
#    (removes the trailing " ENTRY 2" run entirely)
$d.Content.Find.Execute(" ENTRY 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# 2) ", No. 1:17-cv-1519 (DC Cir. 2023). ENTRY 2" -> ", No. 1:17-cv-1519 (DC Cir. 2023)."
#    (removes the trailing space + "ENTRY 2" run, collapsing to a single trailing period)
$d.Content.Find.Execute("2023). ENTRY 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023).", 2) | Out-Null
